# Apply updated dSF (column F) values as part of a data repull / recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    4  = 5
    5  = 0
    6  = -1
    7  = -1
    8  = 0
    9  = -2
    10 = 1
    11 = -1
    12 = -3
    13 = 7
    16 = -1
    17 = 4
    21 = -5
    22 = 5
    24 = 4
    25 = 3
    27 = -1
    30 = -5
    31 = -1
    37 = 6
    38 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
